$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $ws.Range('ZZ1').Formula = '="' + $text + '"'
    $ws.Range('ZZ1').Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $ws.Range('ZZ1').ClearContents()
}

Set-TextValue 'D2' '56.499.97'
Set-TextValue 'E2' '  -1.76%  '
Set-TextValue 'D3' '2.381.18'
Set-TextValue 'E3' '  -2.81%  '
Set-TextValue 'E4' '  -0.06%  '
Set-TextValue 'D5' '502.62'
Set-TextValue 'E5' '  -0.92%  '
Set-TextValue 'D6' '130.38'
Set-TextValue 'E6' '  -1.96%  '
Set-TextValue 'E7' '  +0.17%  '
Set-TextValue 'E8' '  -2.30%  '
Set-TextValue 'D9' '2.389.60'
Set-TextValue 'E9' '  -2.39%  '
Set-TextValue 'D10' '0.0990'
Set-TextValue 'E10' '  +0.79%  '
Set-TextValue 'E11' '  +0.69%  '
Set-TextValue 'D12' '0.329'
Set-TextValue 'E12' '  +2.27%  '
Set-TextValue 'D13' '4.68'
Set-TextValue 'E13' '  +0.91%  '
Set-TextValue 'D14' '2.803.46'
Set-TextValue 'E14' '  -1.51%  '
Set-TextValue 'D15' '56.461.12'
Set-TextValue 'E15' '  -1.83%  '
Set-TextValue 'D16' '21.71'
Set-TextValue 'E16' '  -1.23%  '
Set-TextValue 'E17' '  -0.82%  '
Set-TextValue 'D18' '2.309.59'
Set-TextValue 'E18' '  -4.65%  '
Set-TextValue 'E19' '  -2.48%  '
Set-TextValue 'D20' '4.03'
Set-TextValue 'E20' '  -2.25%  '
Set-TextValue 'D21' '308.01'
Set-TextValue 'E21' '  -2.11%  '
Set-TextValue 'D22' '6.27'
Set-TextValue 'E22' '  -1.78%  '
Set-TextValue 'E23' '  +0.10%  '
Set-TextValue 'D24' '65.25'
Set-TextValue 'E24' '  -0.23%  '
Set-TextValue 'E25' '  +0.32%  '
Set-TextValue 'D26' '0.368'
Set-TextValue 'E26' '  -3.68%  '
Set-TextValue 'D27' '0.148'
Set-TextValue 'E27' '  -3.83%  '
Set-TextValue 'D28' '7.32'
Set-TextValue 'E28' '  -3.96%  '
Set-TextValue 'D29' '172.73'
Set-TextValue 'E29' '  -0.64%  '
Set-TextValue 'D30' '0.0₃0717'
Set-TextValue 'E30' '  -2.60%  '
Set-TextValue 'E31' '  -2.94%  '
Set-TextValue 'E32' '  +0.25%  '
Set-TextValue 'D33' '5.78'
Set-TextValue 'E33' '  -6.92%  '
Set-TextValue 'E34' '  -4.42%  '
Set-TextValue 'E35' '  +0.06%  '
Set-TextValue 'D36' '17.63'
Set-TextValue 'E36' '  -1.93%  '
Set-TextValue 'E37' '  -5.87%  '
Set-TextValue 'D38' '3.79'
Set-TextValue 'E38' '  -2.65%  '
Set-TextValue 'D39' '36.08'
Set-TextValue 'E39' '  -1.17%  '
Set-TextValue 'D40' '0.796'
Set-TextValue 'E40' '  -3.27%  '
Set-TextValue 'E41' '  -4.54%  '
Set-TextValue 'D42' '131.17'
Set-TextValue 'E42' '  -2.94%  '
Set-TextValue 'D43' '3.37'
Set-TextValue 'E43' '  -1.15%  '
Set-TextValue 'E44' '  -4.49%  '
Set-TextValue 'E45' '  -0.67%  '
Set-TextValue 'D46' '0.0909'
Set-TextValue 'E46' '  -0.99%  '
Set-TextValue 'D47' '241.98'
Set-TextValue 'E47' '  -5.93%  '
Set-TextValue 'D48' '0.0485'
Set-TextValue 'E48' '  -2.09%  '
Set-TextValue 'E49' '  -2.01%  '
Set-TextValue 'D50' '17.19'
Set-TextValue 'E50' '  +0.39%  '
Set-TextValue 'E51' '  -3.02%  '
